$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "1.00", "6.20")
# need to be forced to Text format first, otherwise Excel auto-converts the
# string into a numeric value and loses the trailing/structural zeros.
# NumberFormat is applied per-cell (not via a unioned Range) because the
# unioned Range(".. , ..") form only formats the first area.
$textForceRefs = @("D4", "D5", "D6", "D9", "D11", "D13", "D18", "D19", "D21", "D23", "D24", "D25", "D27", "D28", "D30", "D31", "D32", "D33", "D37", "D38", "D39", "D41", "D42", "D43", "D45", "D46", "D48", "D49", "D51")
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.441.70"
$ws.Range("E2").Value = "  -3.30%  "
$ws.Range("D3").Value = "3.702.63"
$ws.Range("E3").Value = "  -3.63%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "595.93"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").Value = "165.41"
$ws.Range("E6").Value = "  -5.14%  "
$ws.Range("D7").Value = "3.704.05"
$ws.Range("E7").Value = "  -3.56%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("D11").Value = "6.20"
$ws.Range("E11").Value = "  -3.90%  "
$ws.Range("E12").Value = "  -3.54%  "
$ws.Range("D13").Value = "37.57"
$ws.Range("E13").Value = "  -5.75%  "
$ws.Range("E14").Value = "  -5.11%  "
$ws.Range("D15").Value = "4.320.01"
$ws.Range("E15").Value = "  -3.57%  "
$ws.Range("D16").Value = "3.697.94"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").Value = "67.490.04"
$ws.Range("E17").Value = "  -3.33%  "
$ws.Range("D18").Value = "17.59"
$ws.Range("E18").Value = "  +5.96%  "
$ws.Range("D19").Value = "7.17"
$ws.Range("E19").Value = "  -3.75%  "
$ws.Range("E20").Value = "  -3.33%  "
$ws.Range("D21").Value = "491.55"
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("E22").Value = "  -3.30%  "
$ws.Range("D23").Value = "0.725"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").Value = "85.91"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "2.30"
$ws.Range("E25").Value = "  -5.97%  "
$ws.Range("E26").Value = "  -2.02%  "
$ws.Range("D27").Value = "12.17"
$ws.Range("E27").Value = "  -3.47%  "
$ws.Range("D28").Value = "10.15"
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "2.93"
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("D31").Value = "2.35"
$ws.Range("E31").Value = "  -6.43%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "31.64"
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "7.63"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("D34").Value = "3.839.52"
$ws.Range("E35").Value = "  -4.48%  "
$ws.Range("D36").Value = "3.642.04"
$ws.Range("E36").Value = "  -3.57%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -3.76%  "
$ws.Range("D39").Value = "5.75"
$ws.Range("E39").Value = "  -5.16%  "
$ws.Range("E40").Value = "  -6.76%  "
$ws.Range("D41").Value = "0.322"
$ws.Range("E41").Value = "  -3.80%  "
$ws.Range("D42").Value = "434.72"
$ws.Range("E42").Value = "  -10.63%  "
$ws.Range("D43").Value = "48.66"
$ws.Range("E44").Value = "  -5.86%  "
$ws.Range("D45").Value = "2.78"
$ws.Range("E45").Value = "  -6.50%  "
$ws.Range("D46").Value = "8.38"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "40.54"
$ws.Range("E48").Value = "  -5.92%  "
$ws.Range("D49").Value = "142.76"
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("D50").Value = "2.757.86"
$ws.Range("E50").Value = "  -5.69%  "
$ws.Range("D51").Value = "0.0347"
$ws.Range("E51").Value = "  -3.58%  "
